$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the line-number references in the "代码位置" (B) column for rows 31-35
# (Verilog source lines shifted up by 14 after an earlier edit to the .v file)
$ws.Range("B31").Value = "ysyx_210611.v:2782:35: Bits of signal are not used: 'itrp_info'[11:8,6:0]"
$ws.Range("B32").Value = "ysyx_210611.v:2814:63: Bits of signal are not used: 'mie_rd_data'[63:12,10:8,6:4,2:0]"
$ws.Range("B33").Value = "ysyx_210611.v:2814:76: Bits of signal are not used: 'mip_rd_data'[63:12,10:8,6:4,2:0]"
$ws.Range("B34").Value = "ysyx_210611.v:2947:16: Bits of signal are not used: 'mem_reg_wr_ctrl'[2]"
$ws.Range("B35").Value = "ysyx_210611.v:2947:33: Bits of signal are not used: 'wb_reg_wr_ctrl'[2]"

# Update the saved scroll position / selection shown when the sheet is reopened:
# viewport top-left moves from A25 to A22, and the selected cell moves from B35 to B30.
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B30").Select()
